$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on price cells whose new values are purely
# numeric-looking, so Excel keeps storing them as text (matching the
# original inline-string cell type) instead of silently coercing them
# into floating point numbers (e.g. '4.90' -> 4.9, '1.00' -> 1).
$textCells = @("D5", "D6", "D8", "D10", "D11", "D12", "D13", "D16", "D17", "D19", "D21", "D22", "D23", "D25", "D26", "D28", "D31", "D32", "D34", "D36", "D38", "D39", "D42", "D43", "D44", "D46", "D47")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Update coin table (rows 2-51) with latest price / volume(1h) figures,
# and swap the RenderToken / Stacks rows (42 <-> 43) to match new ranking.
$ws.Range("D2").Value = '69.410.29'
$ws.Range("E2").Value = '  -1.27%  '
$ws.Range("D3").Value = '2.490.06'
$ws.Range("E3").Value = '  -1.31%  '
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").Value = '565.78'
$ws.Range("E5").Value = '  -1.40%  '
$ws.Range("D6").Value = '163.72'
$ws.Range("E6").Value = '  -2.57%  '
$ws.Range("E7").Value = '  -0.20%  '
$ws.Range("D8").Value = '0.511'
$ws.Range("E8").Value = '  -1.98%  '
$ws.Range("D9").Value = '2.488.10'
$ws.Range("E9").Value = '  -1.15%  '
$ws.Range("D10").Value = '0.158'
$ws.Range("E10").Value = '  -3.24%  '
$ws.Range("D11").Value = '0.166'
$ws.Range("E11").Value = '  -0.59%  '
$ws.Range("D12").Value = '0.353'
$ws.Range("E12").Value = '  -0.95%  '
$ws.Range("D13").Value = '4.90'
$ws.Range("E13").Value = '  -0.47%  '
$ws.Range("D14").Value = '2.943.60'
$ws.Range("E14").Value = '  -1.97%  '
$ws.Range("D15").Value = '69.263.10'
$ws.Range("E15").Value = '  -1.71%  '
$ws.Range("D16").Value = '0.0000175'
$ws.Range("E16").Value = '  -2.07%  '
$ws.Range("D17").Value = '24.25'
$ws.Range("E17").Value = '  -3.73%  '
$ws.Range("D18").Value = '2.468.17'
$ws.Range("E18").Value = '  -3.16%  '
$ws.Range("D19").Value = '11.18'
$ws.Range("E19").Value = '  -2.31%  '
$ws.Range("E20").Value = '  -6.98%  '
$ws.Range("D21").Value = '344.28'
$ws.Range("E21").Value = '  -2.02%  '
$ws.Range("D22").Value = '3.86'
$ws.Range("E22").Value = '  -1.71%  '
$ws.Range("D23").Value = '1.92'
$ws.Range("E23").Value = '  -4.31%  '
$ws.Range("E24").Value = '  -0.17%  '
$ws.Range("D25").Value = '69.47'
$ws.Range("E25").Value = '  -1.31%  '
$ws.Range("D26").Value = '3.88'
$ws.Range("E26").Value = '  -3.38%  '
$ws.Range("D27").Value = '2.616.03'
$ws.Range("E27").Value = '  -3.15%  '
$ws.Range("D28").Value = '8.65'
$ws.Range("E28").Value = '  -2.84%  '
$ws.Range("E29").Value = '  -0.22%  '
$ws.Range("D30").Value = '0.0₃0872'
$ws.Range("E30").Value = '  -4.15%  '
$ws.Range("D31").Value = '7.67'
$ws.Range("E31").Value = '  -2.95%  '
$ws.Range("D32").Value = '441.73'
$ws.Range("E32").Value = '  -5.07%  '
$ws.Range("E33").Value = '  -6.15%  '
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  -0.26%  '
$ws.Range("E35").Value = '  -2.62%  '
$ws.Range("D36").Value = '156.10'
$ws.Range("E36").Value = '  -1.27%  '
$ws.Range("E37").Value = '  -4.35%  '
$ws.Range("D38").Value = '19.05'
$ws.Range("E38").Value = '  -0.22%  '
$ws.Range("D39").Value = '18.12'
$ws.Range("E39").Value = '  -3.38%  '
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("E41").Value = '  -1.50%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D42").Value = '4.58'
$ws.Range("E42").Value = '  -4.79%  '
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").Value = '1.58'
$ws.Range("E43").Value = '  -1.57%  '
$ws.Range("D44").Value = '2.15'
$ws.Range("E44").Value = '  -7.50%  '
$ws.Range("E45").Value = '  -7.78%  '
$ws.Range("D46").Value = '138.28'
$ws.Range("E46").Value = '  -3.35%  '
$ws.Range("D47").Value = '3.43'
$ws.Range("E47").Value = '  -2.24%  '
$ws.Range("E48").Value = '  -3.55%  '
$ws.Range("E49").Value = '  -0.93%  '
$ws.Range("E50").Value = '  -1.17%  '
$ws.Range("E51").Value = '  -1.27%  '
